$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new parameter row ("EmailRemetentes") after the "EmailPasta" row.
# This shifts the previous rows 6-14 down to 7-15, carrying their formatting.
$ws.Rows.Item(6).Insert()

$ws.Cells.Item(6, 1).Value2 = "EmailRemetentes"
$ws.Cells.Item(6, 3).Value2 = "Lista de Remetentes que enviam o e-mail desejado, separados por ;"
$ws.Rows.Item(6).RowHeight = 14.25

# Rename a few existing parameter names for consistency ("Otimizados parametros de configuracao").
$ws.Cells.Item(3, 1).Value2 = "CaminhoOutputs"
$ws.Cells.Item(9, 1).Value2 = "AWSLogin"
$ws.Cells.Item(11, 1).Value2 = "AWSBucketName"

# Widen column C to fit the new, longer description text.
$ws.Columns.Item(3).ColumnWidth = 64.5

# Leave the selection where the edit was made.
$ws.Range("C10").Select()
